$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill B2:K2 with "AFM"
$ws.Range("B2:K2").Value = "AFM"

# Fill B3:D3 with "AFM"
$ws.Range("B3:D3").Value = "AFM"

# Update the active cell / selection to G26
$ws.Range("G26").Select()
